$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.213415622711182
$ws.Range("B1").Value = 2.058534383773804
$ws.Range("C1").Value = 4.429384231567383
$ws.Range("D1").Value = 2.95613169670105
$ws.Range("E1").Value = 1.194478034973145
